$d = $word.ActiveDocument

# Update the header date (unique text in the document, safe to use document-wide Find/Replace)
$d.Content.Find.Execute("2026-02-24 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-25 Wednesday", 2)

$t = $d.Tables.Item(1)

# Helper: replace the old expression with the new one inside a single table cell,
# without relying on Find.Execute (which searches the whole document rather than
# being scoped to the calling Range in this runtime). Operating on Range.Text keeps
# the existing run formatting (font/size) untouched.
function Replace-CellText($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.Text = $r.Text.Replace($oldText, $newText)
}

# Row 1
Replace-CellText $t 1 1 "237×3=" "544×4="
Replace-CellText $t 1 2 "502×9=" "476×4="
Replace-CellText $t 1 3 "147×6=" "766×2="
Replace-CellText $t 1 4 "311×2=" "559×3="
Replace-CellText $t 1 5 "646×7=" "568×7="

# Row 5
Replace-CellText $t 5 1 "888×8=" "855×5="
Replace-CellText $t 5 2 "177×7=" "817×2="
Replace-CellText $t 5 3 "276×4=" "835×3="
Replace-CellText $t 5 4 "126×3=" "263×2="
Replace-CellText $t 5 5 "494×5=" "754×6="

# Row 10
Replace-CellText $t 10 1 "365×7=" "687×4="
Replace-CellText $t 10 2 "636×8=" "119×2="
Replace-CellText $t 10 3 "138×8=" "850×4="
Replace-CellText $t 10 4 "772×5=" "413×3="
Replace-CellText $t 10 5 "308×8=" "508×6="

# Row 15
Replace-CellText $t 15 1 "905×2=" "897×2="
Replace-CellText $t 15 2 "367×6=" "614×5="
Replace-CellText $t 15 3 "246×7=" "288×2="
Replace-CellText $t 15 4 "413×3=" "565×2="
Replace-CellText $t 15 5 "875×7=" "929×8="

# Row 20
Replace-CellText $t 20 1 "596×9=" "471×8="
Replace-CellText $t 20 2 "693×5=" "954×3="
Replace-CellText $t 20 3 "311×9=" "652×7="
Replace-CellText $t 20 4 "938×9=" "812×3="
Replace-CellText $t 20 5 "618×2=" "636×7="
